$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtergebnis")

# Make the overall-results sheet ("Gesamtergebnis") the active sheet/tab.
$ws.Activate()

# The balance columns (Startguthaben / Endsaldo) on the "Total" row were
# previously left as "N/A" text; now include real totals (0) for them too.
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# Move the active selection to C5 on this sheet.
$ws.Range("C5").Select()
